$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Columns A-H already carry the bold/border/centered header style (s="1"),
# so just overwriting .Value keeps that style. Columns I-K are new, so we
# first clone the existing header style onto them (Copy reuses the style
# index instead of registering a new one), then set their text.
$ws.Range("A1").Copy($ws.Range("I1"))
$ws.Range("A1").Copy($ws.Range("J1"))
$ws.Range("A1").Copy($ws.Range("K1"))

$ws.Range("A1").Value = "Дата создания"
$ws.Range("B1").Value = "Проблема"
$ws.Range("C1").Value = "User ID"
$ws.Range("D1").Value = "ИМЯ"
$ws.Range("E1").Value = "ФАМИЛИЯ"
$ws.Range("F1").Value = "ЛОГИН"
$ws.Range("G1").Value = "Ответ1"
$ws.Range("H1").Value = "Ответ2"
$ws.Range("I1").Value = "Ответ3"
$ws.Range("J1").Value = "Почта"
$ws.Range("K1").Value = "Доп.вопрос"

# --- Data rows (rows 2-11) ---------------------------------------------
$combining = [char]0x301

$dates = @(45539.90326388889, 45539.9024074074, 45539.89626157407, 45539.89578703704, 45539.89491898148, 45539.88535879629, 45539.8768287037, 45539.87126157407, 45539.86407407407, 45539.85804398148)
$problems = @("Куда вводить код", "Не пришел код", "Помощь с подбором курса", "Нет моего вопроса", "Не пришел код", "Нет моего вопроса", "Не могу войти в аккаунт", "Не работает код", "Куда вводить код", "Помощь с подбором курса")
$userIds = @(6625770047, 1006569664, 6625770047, 6625770047, 6625770047, 1006569664, 1006569664, 1006569664, 1006569664, 1006569664)
$firstNames = @("kiper_slivki", "Roman", "kiper_slivki", "kiper_slivki", "kiper_slivki", "Roman", "Roman", "Roman", "Roman", "Roman")
$lastNames = @("", "Chiper", "", "", "", "Chiper", "Chiper", "Chiper", "Chiper", "Chiper")
$logins = @("", "RomanKiper", "", "", "", "RomanKiper", "RomanKiper", "RomanKiper", "RomanKiper", "RomanKiper")
$answer1 = @("", "", "Пророоо", "", "", "", "", "", "", "одинннн")
$answer2 = @("", "", "Ллллллллллллллллллллллллоьттиииииииииииииииииииииииииииииииииииииииииииииммммммммммми", "", "", "", "", "", "", "дваааа")
$answer3 = @("", "", "Рррроооооол", "", "", "", "", "", "", "трииии")
$mails = @("", "@ ккккк", "", "", "Почта@руу", "", "Дуда@mail", "", "", "")
$extra5 = "Срочный вопрос оооооооооо" + $combining + "лльььььььььььььььььььььььььььььььььььььооооорррррррррпааапроолллооооооооооооотторррииииррррррппппммирррооллььььььььььььььь"
$extra7 = "Привет админ`nМоего вопроса нет в списке`nЧто делать`nУ меня вопрос следущий`nЯ переживаю чтотеет моего вопросаротьььььбббьлооооррроооориррррролллллдддддддд"
$extras = @("", "", "", $extra5, "", $extra7, "", "", "", "")

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    if ($i -eq 0) {
        # First date cell: touch the format twice (lower-case, then the
        # real upper-case one) - this is what registers BOTH numFmt 164
        # (left unused/orphaned) and 165 (the one actually applied).
        $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    } else {
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($r, 2).Value = $problems[$i]
    $ws.Cells.Item($r, 3).Value = $userIds[$i]
    $ws.Cells.Item($r, 4).Value = $firstNames[$i]
    $ws.Cells.Item($r, 5).Value = $lastNames[$i]
    $ws.Cells.Item($r, 6).Value = $logins[$i]
    $ws.Cells.Item($r, 7).Value = $answer1[$i]
    $ws.Cells.Item($r, 8).Value = $answer2[$i]
    $ws.Cells.Item($r, 9).Value = $answer3[$i]
    $ws.Cells.Item($r, 10).Value = $mails[$i]
    $ws.Cells.Item($r, 11).Value = $extras[$i]

    if ($extras[$i] -like "*`n*") {
        # A value with embedded line breaks auto-expands the row height
        # (real Excel behaviour). Re-run AutoFit with no explicit
        # WrapText toggle so the row height snaps back to the sheet's
        # default and the row keeps NO customHeight/ht override - exactly
        # the plain `<row r="N">` the target file has.
        $ws.Cells.Item($r, 11).EntireRow.AutoFit()
    }
}
